$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 14725
$ws1.Range("F8").Value = 401
$ws1.Range("F10").Value = 15183
$ws1.Range("F12").Value = 8630
$ws1.Range("F16").Value = 178
$ws1.Range("F20").Value = 11
$ws1.Range("F22").Value = 21
$ws1.Range("F25").Value = 1072
$ws1.Range("F28").Value = 53
$ws1.Range("F30").Value = 29
$ws1.Range("F31").Value = 417
$ws1.Range("F32").Value = 24
$ws1.Range("F33").Value = 28
$ws1.Range("F35").Value = 265
$ws1.Range("F36").Value = 419
$ws1.Range("F38").Value = 5328
$ws1.Range("F39").Value = 5226

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 59

# --- Sheet "全部类型" (all types, aggregated) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 14727
$ws4.Range("F8").Value = 401
$ws4.Range("F10").Value = 15183
$ws4.Range("F12").Value = 8630
$ws4.Range("F17").Value = 178
$ws4.Range("F21").Value = 11
$ws4.Range("F23").Value = 21
$ws4.Range("F26").Value = 1072
$ws4.Range("F29").Value = 53
$ws4.Range("F31").Value = 29
$ws4.Range("F32").Value = 59
$ws4.Range("F34").Value = 417
$ws4.Range("F35").Value = 24
$ws4.Range("F36").Value = 28
$ws4.Range("F38").Value = 265
$ws4.Range("F39").Value = 419
$ws4.Range("F41").Value = 5328
$ws4.Range("F42").Value = 5226
